# Updates the cryptos price list (Sheet1) with refreshed Price/Volume(1h) values,
# and swaps the Fetch.AI / RenzoRestakedETH rows (33 & 34) back to their
# previous ranking order, as captured by the latest GitHub Actions data pull.
#
# For numeric-looking Price strings (column D) we force the cell's number
# format to Text ("@") before assigning the value, so Excel stores the exact
# source string (e.g. "125.00", "1.38") instead of silently coercing it to a
# number and dropping trailing zeros / decimal formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.389.84"
$ws.Range("E2").Value = "  -3.12%  "
$ws.Range("D3").Value = "3.378.78"
$ws.Range("E3").Value = "  -3.76%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.26"
$ws.Range("E5").Value = "  -3.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.00"
$ws.Range("E6").Value = "  -7.09%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.380.47"
$ws.Range("E8").Value = "  -3.69%  "
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("E10").Value = "  -4.88%  "
$ws.Range("E11").Value = "  -4.41%  "
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("D13").Value = "3.951.75"
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").Value = "3.378.45"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D17").Value = "62.439.98"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.41"
$ws.Range("E18").Value = "  -5.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.17"
$ws.Range("E19").Value = "  -8.91%  "
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.51"
$ws.Range("E23").Value = "  -4.34%  "
$ws.Range("D24").Value = "3.515.63"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.13"
$ws.Range("E26").Value = "  -4.81%  "
$ws.Range("E27").Value = "  -10.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.88"
$ws.Range("E29").Value = "  -7.00%  "
$ws.Range("E30").Value = "  -6.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.74"
$ws.Range("E31").Value = "  -6.41%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.38"
$ws.Range("E33").Value = "  -5.59%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.409.91"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("E35").Value = "  -6.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.60"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.19"
$ws.Range("E37").Value = "  -3.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "164.93"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.59"
$ws.Range("E39").Value = "  -5.39%  "
$ws.Range("E40").Value = "  -5.66%  "
$ws.Range("E41").Value = "  -4.96%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -5.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.41"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.23"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.54"
$ws.Range("E46").Value = "  -7.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.44"
$ws.Range("E47").Value = "  -9.88%  "
$ws.Range("E48").Value = "  -8.79%  "
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").Value = "2.233.96"
$ws.Range("E50").Value = "  -6.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.838"
$ws.Range("E51").Value = "  -6.74%  "
